$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.806.05"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "'1.611.16"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'210.89"
$ws.Range("E5").Value = "  -2.48%  "

$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("D9").Value = "'0.0621"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "'19.71"
$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("E11").Value = "  -1.34%  "

$ws.Range("D12").Value = "'1.838.78"
$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").Value = "'1.602.01"
$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'26.792.63"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'63.86"
$ws.Range("E17").Value = "  -2.96%  "

$ws.Range("D18").Value = "'0.0₃0730"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").Value = "'210.57"
$ws.Range("E19").Value = "  -1.57%  "

$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").Value = "'2.33"
$ws.Range("E23").Value = "  -6.93%  "

$ws.Range("D24").Value = "'8.88"
$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("D25").Value = "'146.66"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").Value = "'7.48"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("D29").Value = "'15.37"
$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  -2.47%  "

$ws.Range("D33").Value = "'0.687"
$ws.Range("E33").Value = "  +26.87%  "

$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("D35").Value = "'1.319.77"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  -1.47%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D41").Value = "'0.793"
$ws.Range("E41").Value = "  -1.72%  "

$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("D43").Value = "'5.30"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").Value = "'63.32"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("D45").Value = "'1.749.26"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").Value = "'89.04"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("D48").Value = "'0.817"
$ws.Range("E48").Value = "  +5.74%  "

$ws.Range("D49").Value = "'0.0511"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("D50").Value = "'0.0977"
$ws.Range("E50").Value = "  +3.24%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.0₇0980"
$ws.Range("E51").Value = "  -3.12%  "
